# Backup QR Scanner data - append newly scanned log rows to the
# Pathology Lab/Museum sheet (rows 168 and 169).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scan log entries, in column order: Student ID, Subject, Log Date,
# Log Time, Type, User.
$newRows = @(
    @("244111", "Pathology Lab/Museum", "18/11/2025", "11:11:44", "Scan", "mona.I.hussein@med.asu.edu.eg"),
    @("244055", "Pathology Lab/Museum", "18/11/2025", "11:11:52", "Scan", "mona.I.hussein@med.asu.edu.eg")
)

# Find the first empty row right after the existing data (row 167 in the
# source workbook) and append the new rows there.
$lastUsedRow = $ws.UsedRange.Rows.Count
if ($lastUsedRow -lt 167) {
    $lastUsedRow = 167
}
$startRow = $lastUsedRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]

    # Column A holds numeric-looking IDs that must stay stored as text
    # (matching the rest of the column), so force a text number format
    # before assigning the value - otherwise Excel would coerce it to a
    # number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $values[0]

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
